# feat: add 2022-Q4 data
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update "总计" (summary) sheet: shift the existing 2021-Q4 / 2021-Q3
#    rows down by one (row3 <- old row2, row4 <- old row3) and put the
#    new 2022-Q4 figures into row 2. Writing the cells directly (rather
#    than Rows.Insert()) avoids picking up stray inherited formatting.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

# Row 4 (new) needs the same index-column style as row 3 / row 2 (A3).
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A4").PasteSpecial(-4122)

$wsTotal.Range("A4").Value = 2
$wsTotal.Range("B4").Value = "2021-Q3"
$wsTotal.Range("C4").Value = 7
$wsTotal.Range("D4").Value = 0.05

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2021-Q4"
$wsTotal.Range("C3").Value = 2
$wsTotal.Range("D3").Value = 0.01

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 1
$wsTotal.Range("D2").Value = 0.06

# ---------------------------------------------------------------------
# 2) Insert a brand-new "2022-Q4" worksheet right before the existing
#    "2021-Q4" sheet.
# ---------------------------------------------------------------------
$wsOldQ4Ref = $wb.Worksheets.Item("2021-Q4")
$wsNew = $wb.Worksheets.Add($wsOldQ4Ref)
$wsNew.Name = "2022-Q4"

# Re-fetch a fresh, stable reference to the (now shifted) 2021-Q4 sheet.
$wsOldQ4 = $wb.Worksheets.Item("2021-Q4")

# Header row (B1:H1) - copy style from the existing 2021-Q4 sheet header
$wsOldQ4.Range("B1:H1").Copy()
$wsNew.Range("B1:H1").PasteSpecial(-4122)

$wsNew.Range("B1").Value = "基金代码"
$wsNew.Range("C1").Value = "基金名称"
$wsNew.Range("D1").Value = "基金规模"
$wsNew.Range("E1").Value = "股票总仓位"
$wsNew.Range("F1").Value = "仓位占比"
$wsNew.Range("G1").Value = "持有市值(亿元)"
$wsNew.Range("H1").Value = "仓位排名"

# Data row 2 - A2 style matches the existing 2021-Q4 sheet's A2 style
$wsOldQ4.Range("A2").Copy()
$wsNew.Range("A2").PasteSpecial(-4122)
$wsNew.Range("A2").Value = 0

# Text-ish numeric columns must stay text (preserve "000974", "2.33", ...)
$wsNew.Range("B2:G2").NumberFormat = "@"
$wsNew.Range("B2").Value = "000974"
$wsNew.Range("C2").Value = "安信消费医药主题股票"
$wsNew.Range("D2").Value = "2.33"
$wsNew.Range("E2").Value = "92.45"
$wsNew.Range("F2").Value = "2.53"
$wsNew.Range("G2").Value = "0.0589"
$wsNew.Range("B2:G2").ClearFormats()
$wsNew.Range("H2").Value = 9

# Restore the originally-active tab (2021-Q3), since adding the new
# sheet shifts Excel's active-tab selection to it by default.
$wsOldQ3 = $wb.Worksheets.Item("2021-Q3")
$wsOldQ3.Activate()

Write-Host "done"
